# Updating hands-on with last package versions.
#
# 1) Insert a new "Thanks to..." slide (Title and Content layout) as the
#    2nd slide of the deck (it becomes sldId 262, positioned right after
#    the title slide and before the two "Machine Learning Concepts
#    Taxonomy" slides).
# 2) Bump the fixed header/footer date from 11.01.2024 to 12.01.2024 on
#    the slide master and every slide layout.

$p = $ppt.ActivePresentation

# --- 1) New slide -----------------------------------------------------
$s = $p.Slides.Add(2, 2)   # 2 = ppLayoutText -> "Title and Content" custom layout

# Title placeholder
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Thanks to…"

# Body / content placeholder
$body = $s.Shapes.Item(2).TextFrame
$body.TextRange.Text = "Organisers`rSusanna Bisogni`rAdriana Gargiulo`rMarco Scodeggio`rMarco Fumana`r`rAnd to the entire INAF for the support."

# Paragraph 1 ("Organisers") is a plain (non-bulleted) lead-in line.
$body.TextRange.Paragraphs(1,1).ParagraphFormat.Bullet.Type = 0
# Paragraph 6 is the blank spacer line.
$body.TextRange.Paragraphs(6,1).ParagraphFormat.Bullet.Type = 0
# Paragraph 7 ("And to the entire INAF for the support.") is also plain.
$body.TextRange.Paragraphs(7,1).ParagraphFormat.Bullet.Type = 0

# --- 2) Fixed date placeholders: 11.01.2024 -> 12.01.2024 -------------
$newDate = "12.01.2024"

$master = $p.SlideMaster
for ($j = 1; $j -le $master.Shapes.Count; $j++) {
    $sh = $master.Shapes.Item($j)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $cl = $master.CustomLayouts.Item($i)
    for ($j = 1; $j -le $cl.Shapes.Count; $j++) {
        $sh = $cl.Shapes.Item($j)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}
